$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate row 2's data into a new row 3 (same values/types as row 2).
$ws.Range("A2:Q2").Copy()
$ws.Range("A3:Q3").PasteSpecial()
